$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "1OHXOK"
$ws.Range("B22").Value = "Cinta Flex Epson"
$ws.Range("C22").Value = "TM U220"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 50000
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 12
$ws.Range("H22").Formula = "=(E22-D22)*G22"
$ws.Range("I22").Formula = "=D22*F22"
$ws.Range("J22").Value = 0
